# Daily task.xlsx update:
#  - append new rows to the "python" sheet tracking progress on "loops" /
#    Dataiku videos
#  - add a brand-new "Tablaeu" worksheet (after "python") with its own
#    Date/Topic/Status tracker table

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "python"

# ---------------------------------------------------------------------
# 1. "python" sheet - add rows 5, 6 and 7 under the existing table
# ---------------------------------------------------------------------

# Row 5: 1/9/2025 - Dataiku vidoes 1 - completed
$ws1.Range("A5").Value2 = 45666
$ws1.Range("A5").NumberFormat = "m/d/yy"
$ws1.Range("B5").Value = "Dataiku vidoes 1"
$ws1.Range("D5").Value = "completed"

# Row 6: 1/10/2025 - Dataiku vidoes 2 - completed (date cell vertically centered)
$ws1.Range("A6").Value2 = 45667
$ws1.Range("A6").NumberFormat = "m/d/yy"
$ws1.Range("A6").VerticalAlignment = -4108  # xlCenter
$ws1.Range("B6").Value = "Dataiku vidoes 2"
$ws1.Range("D6").Value = "completed"

# Row 7: 1/10/2025 - loops - completed (date cell vertically centered)
$ws1.Range("A7").Value2 = 45667
$ws1.Range("A7").NumberFormat = "m/d/yy"
$ws1.Range("A7").VerticalAlignment = -4108  # xlCenter
$ws1.Range("B7").Value = "loops"
$ws1.Range("D7").Value = "completed"

$ws1.Range("A8").Select()

# ---------------------------------------------------------------------
# 2. Add the new "Tablaeu" worksheet right after "python"
# ---------------------------------------------------------------------

$ws2 = $wb.Worksheets.Add([System.Type]::Missing, $ws1)
$ws2.Name = "Tablaeu"

$ws2.Columns.Item(1).ColumnWidth = 13.333333333333334
$ws2.Columns.Item(2).ColumnWidth = 25
$ws2.Columns.Item(3).ColumnWidth = 25.666666666666668
$ws2.Columns.Item(4).ColumnWidth = 12.166666666666666

# Header row
$ws2.Range("A1").Value = "Date"
$ws2.Range("A1").NumberFormat = "m/d/yy"
$ws2.Range("B1").Value = "Topic"
$ws2.Range("C1").Value = "Status"
$ws2.Range("C1").HorizontalAlignment = -4152  # xlRight
$ws2.Range("C1").VerticalAlignment = -4108    # xlCenter

# Row 2: 10/7/2025 - tableau introduction - inprogress - completed
$ws2.Range("A2").Value2 = 45937
$ws2.Range("A2").NumberFormat = "m/d/yy"
$ws2.Range("B2").Value = "tableau introduction"
$ws2.Range("C2").Value = "inprogress"
$ws2.Range("D2").Value = "completed"

# Row 3: 10/8/2025 - tableau dashboatd creation - completed
$ws2.Range("A3").Value2 = 45938
$ws2.Range("A3").NumberFormat = "m/d/yy"
$ws2.Range("B3").Value = "tableau dashboatd creation"
$ws2.Range("D3").Value = "completed"

# Row 4: 1/9/2025 - tableau filters - completed
$ws2.Range("A4").Value2 = 45666
$ws2.Range("A4").NumberFormat = "m/d/yy"
$ws2.Range("B4").Value = "tableau filters"
$ws2.Range("D4").Value = "completed"

# Row 5: 1/10/2025 - tablaue jons - completed
$ws2.Range("A5").Value2 = 45667
$ws2.Range("A5").NumberFormat = "m/d/yy"
$ws2.Range("B5").Value = "tablaue jons"
$ws2.Range("D5").Value = "completed"

$ws2.Range("C6").Select()

# Re-activate "python" so it stays the visible/selected tab
$ws1.Activate()
$ws1.Range("A8").Select()
